# Apply the PC_Users login/role naming update.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PC_Users")

# Update LoginId (column C) values to the new Partner Center agent account naming scheme.
$ws.Range("C2").Value = "AdmAgentGlobalAdm"
$ws.Range("C3").Value = "AdmAgentAccAdm"
$ws.Range("C4").Value = "AdmAgentBillAdm"
$ws.Range("C5").Value = "HelpDeskAgent"
$ws.Range("C6").Value = "SalesAgent"

# Row 5 (HelpDesk) LastName also renamed from "Agent" to "HelpDeskAgent".
$ws.Range("B5").Value = "HelpDeskAgent"

# Widen column B so the longer LastName values are fully visible (no longer auto-fit).
$ws.Columns.Item(2).ColumnWidth = 17.1796875

# Move the active selection to C6 (matches the saved selection in the file).
$ws.Range("C6").Select()
